$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")

$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3

$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 2

$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 4

$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2

$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 5

$ws.Range("I4").Select()
